$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.909.69"
$ws.Range("E2").Value = "  +2.60%  "

$ws.Range("D3").Value = "3.594.33"
$ws.Range("E3").Value = "  +2.02%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'600.96"
$ws.Range("E5").Value = "  +2.06%  "

$ws.Range("D6").Value = "'173.94"
$ws.Range("E6").Value = "  +1.26%  "

$ws.Range("D7").Value = "3.587.34"
$ws.Range("E7").Value = "  +1.99%  "

$ws.Range("E8").Value = "  +0.69%  "

$ws.Range("E9").Value = "  -0.02%  "

$ws.Range("D10").Value = "'0.201"
$ws.Range("E10").Value = "  +5.89%  "

$ws.Range("D11").Value = "'7.48"
$ws.Range("E11").Value = "  +8.46%  "

$ws.Range("E12").Value = "  +1.69%  "

$ws.Range("D13").Value = "'47.00"
$ws.Range("E13").Value = "  -0.66%  "

$ws.Range("E14").Value = "  +0.98%  "

$ws.Range("D15").Value = "4.176.24"
$ws.Range("E15").Value = "  +2.25%  "

$ws.Range("D16").Value = "'8.46"
$ws.Range("E16").Value = "  -0.34%  "

$ws.Range("D17").Value = "'615.64"
$ws.Range("E17").Value = "  -1.25%  "

$ws.Range("D18").Value = "3.604.57"
$ws.Range("E18").Value = "  +2.54%  "

$ws.Range("D19").Value = "70.988.07"
$ws.Range("E19").Value = "  +2.64%  "

$ws.Range("E20").Value = "  -1.15%  "

$ws.Range("D21").Value = "'17.53"
$ws.Range("E21").Value = "  +0.58%  "

$ws.Range("E22").Value = "  +0.31%  "

$ws.Range("D23").Value = "'9.27"
$ws.Range("E23").Value = "  -16.99%  "

$ws.Range("D24").Value = "'16.04"
$ws.Range("E24").Value = "  +0.52%  "

$ws.Range("D25").Value = "'97.47"
$ws.Range("E25").Value = "  +0.37%  "

$ws.Range("E26").Value = "  -1.39%  "

$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  -0.04%  "

$ws.Range("E28").Value = "  +0.89%  "

$ws.Range("D29").Value = "'34.12"
$ws.Range("E29").Value = "  +3.98%  "

$ws.Range("D30").Value = "'9.24"
$ws.Range("E30").Value = "  -0.68%  "

$ws.Range("D31").Value = "'8.51"
$ws.Range("E31").Value = "  -0.46%  "

$ws.Range("E32").Value = "  -1.76%  "

$ws.Range("D33").Value = "'7.28"
$ws.Range("E33").Value = "  +4.57%  "

$ws.Range("E34").Value = "  -1.19%  "

$ws.Range("D35").Value = "'633.21"
$ws.Range("E35").Value = "  -0.71%  "

$ws.Range("D36").Value = "'3.74"
$ws.Range("E36").Value = "  +6.61%  "

$ws.Range("E37").Value = "  -1.09%  "

$ws.Range("E38").Value = "  +0.86%  "

$ws.Range("D39").Value = "'0.0481"
$ws.Range("E39").Value = "  +5.53%  "

$ws.Range("D40").Value = "'57.45"
$ws.Range("E40").Value = "  +0.33%  "

$ws.Range("E41").Value = "  +0.09%  "

$ws.Range("E42").Value = "  +5.33%  "

$ws.Range("D43").Value = "3.415.39"
$ws.Range("E43").Value = "  +0.84%  "

$ws.Range("E44").Value = "  -0.80%  "

$ws.Range("D45").Value = "0.0₃0722"
$ws.Range("E45").Value = "  +2.92%  "

$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "'33.11"
$ws.Range("E46").Value = "  +0.63%  "

$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").Value = "'2.98"
$ws.Range("E47").Value = "  +7.96%  "

$ws.Range("D48").Value = "'2.68"
$ws.Range("E48").Value = "  +5.26%  "

$ws.Range("E49").Value = "  +0.64%  "

$ws.Range("D50").Value = "'132.92"
$ws.Range("E50").Value = "  -0.47%  "

$ws.Range("E51").Value = "  -0.08%  "
